$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A23").Value = "Chache on listing"
$ws.Range("A24").Value = "Notification Against Transaction "
$ws.Range("A25").Value = "Notifcication Against Reaction on challenge"
$ws.Range("A26").Value = "Donate within time"
$ws.Range("B26").Value = "Done"
$ws.Range("A27").Value = "post submit challenge"

# Match the author's final selection state (single cell B27)
$ws.Range("B27").Select()
